$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values
$ws.Range("B2").Value = 12.769944043717169
$ws.Range("C2").Value = 11.479533064867338
$ws.Range("D2").Value = 12.668339273527021
$ws.Range("E2").Value = 12.44210363328866

# Row 3 data values
$ws.Range("B3").Value = 12.737997266455428
$ws.Range("C3").Value = 10.896010125357197
$ws.Range("D3").Value = 13.455818205807157
$ws.Range("E3").Value = 11.044606909410176

# Update selection to match the new range
$ws.Range("B1:E3").Select()
